# Rename the RACI header cells to the more descriptive role labels,
# resize the newly-relevant columns, and move the active selection
# to where the author left off (D16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3) relabeling -----------------------------------
# Order matters for shared-string table layout: E3 is written before D3
# so new strings land in the same order as the authoritative edit.
$ws.Range("B3").Value = "MOE (nous)"
$ws.Range("C3").Value = "AMOE (nos employés)"
$ws.Range("E3").Value = "AMOA (quartpi)"
$ws.Range("D3").Value = "MOA (ake michi)"

# --- Column widths for the now-longer header labels -------------------
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666   # B -> ~15.57
$ws.Columns.Item(3).ColumnWidth = 23.5                 # C -> ~24.29
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668   # D -> 19

# --- Restore the author's last active selection ------------------------
[void]$ws.Range("D16").Select()
